# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly scraped values (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> hashtable of cell address -> new value
$updates = @{
    "展览"   = @{ "F2" = 1041; "F3" = 221; "F4" = 2528; "F5" = 40; "F6" = 545 }
    "全部类型" = @{ "F4" = 1041; "F5" = 221; "F6" = 2528; "F7" = 40; "F8" = 545 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($addr in $cellUpdates.Keys) {
        $ws.Range($addr).Value = $cellUpdates[$addr]
    }
}
